# Commit: "removed ER tags from non-ER templates and non-ER tags"
#
# The second worksheet ("SwateTemplateMetadata") is renamed to "isa_template",
# and the "Tags" section of that sheet had a duplicated "ER" tag (PRIDE /
# DPBO_1000098 / DPBO) removed - it was accidentally also listed as a Tag,
# duplicating the dedicated "ER" row above it. Removing it shifts the
# remaining Tag columns (Proteomics/Mass spectrometry/MS/Data/Processing and
# their accession/source-ref columns) one column to the left and drops the
# now-unused trailing column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Rename the metadata sheet.
$ws.Name = "isa_template"

# Row 12 ("Tags"): drop the first tag value (PRIDE, which duplicated the "ER"
# row) and shift the remaining tag values left by one column.
$ws.Range("B12").Value2 = $ws.Range("C12").Value2
$ws.Range("C12").Value2 = $ws.Range("D12").Value2
$ws.Range("D12").Value2 = $ws.Range("E12").Value2
$ws.Range("E12").Value2 = $ws.Range("F12").Value2
$ws.Range("F12").Value2 = $ws.Range("G12").Value2
$ws.Range("G12").ClearContents()

# Row 13 ("Tags Term Accession Number"): drop the accession that belonged to
# the removed PRIDE/DPBO tag, shifting the NCIT accession left into B13.
$ws.Range("B13").Value2 = $ws.Range("C13").Value2
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()

# Row 14 ("Tags Term Source REF"): same shift for the source-ref values.
$ws.Range("B14").Value2 = $ws.Range("C14").Value2
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()

# Column G is now completely unused on this sheet - remove it so the sheet
# dimension shrinks back from A1:G27 to A1:F27.
$ws.Columns("G").Delete()

# Row 12 grew taller (it now wraps across the remaining tag columns), while
# row 13 goes back to the default row height since it lost its second value.
$ws.Rows(12).RowHeight = 72
$ws.Rows(13).AutoFit()

# Update the on-sheet selection to reflect where the author ended up editing.
$ws.Activate()
$ws.Range("B15").Select()
